$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / account holder info ---
$ws.Range("C2").Value = "Hartmut"
# Card number is a long digit string; force text formatting so it isn't
# reinterpreted as a numeric value (matches original inlineStr cell type).
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# --- Opening balance line ---
$ws.Range("D5").Value = "KONTOSTAND AM 12.03.2024"

# --- Row 6 (existing transaction, values change) ---
$ws.Range("B6").Value = "13.03."
$ws.Range("C6").Value = "14.03."
$ws.Range("D6").Value = "BURGER KING Mühldorf am Inn"
$ws.Range("E6").Value = "11,26-"

# --- Row 7 (existing transaction, values change) ---
$ws.Range("B7").Value = "16.03."
$ws.Range("C7").Value = "17.03."
$ws.Range("D7").Value = "RECHNUNG VODAFONE GMBH 26078569"
$ws.Range("E7").Value = "39,61-"

# --- Row 8 (existing transaction, values change) ---
$ws.Range("B8").Value = "18.03."
$ws.Range("C8").Value = "19.03."
$ws.Range("D8").Value = "AMAZON.DE MKTPLC EU VWXGBJ"
$ws.Range("E8").Value = "190,06-"

# --- Rows 9-11 were empty placeholder rows; copy formatting from row 8
#     (which already carries the transaction-row style) before filling values ---
$ws.Range("B8:E8").Copy()
$ws.Range("B9:E9").PasteSpecial(-4122)
$ws.Range("B8:E8").Copy()
$ws.Range("B10:E10").PasteSpecial(-4122)
$ws.Range("B8:E8").Copy()
$ws.Range("B11:E11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B9").Value = "20.03."
$ws.Range("C9").Value = "21.03."
$ws.Range("D9").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E9").Value = "24,99-"

$ws.Range("B10").Value = "24.03."
$ws.Range("C10").Value = "25.03."
$ws.Range("D10").Value = "ZALANDO MKTPLC EU UBEWXL"
$ws.Range("E10").Value = "38,95-"

$ws.Range("B11").Value = "26.03."
$ws.Range("C11").Value = "27.03."
$ws.Range("D11").Value = "PAYPAL CHBHSQ"
$ws.Range("E11").Value = "55,50-"

# --- Closing balance line ---
$ws.Range("D12").Value = "KONTOSTAND AM 29.03.2024"
$ws.Range("E12").Value = "360,37-"

# --- Next billing date note ---
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 08.04.2024"
